$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns: P1=14, Q1=15
# Copy the formatting of O1 (bold header style) into P1:Q1, then set values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing data columns I, K, M, O for rows 2-25, and add new columns P, Q = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1  # column K
    $ws.Cells.Item($r, 13).Value = 2  # column M
    $ws.Cells.Item($r, 15).Value = 1  # column O
    $ws.Cells.Item($r, 16).Value = 2  # column P
    $ws.Cells.Item($r, 17).Value = 2  # column Q
}
